# "Added changable column names"
#
# The header row (row 1) on Arkusz1 used fixed labels WYTW1/WYTW2/WYTW3 for
# columns B:D and had nothing in A1. Relabel them with friendlier / more
# descriptive names and add a caption ("Opis") in A1 describing the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Opis"
$ws.Range("B1").Value = "WytwPierwszy"
$ws.Range("C1").Value = "WytwDrugi"
$ws.Range("D1").Value = "WytwTrzeci"

# Columns B:D were already sized to "best fit" their (now longer/shorter)
# header text - refresh that auto-sizing for the new labels.
$ws.Columns.Item(2).ColumnWidth = 13.3
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 10.5

# Move the active selection, as left by the author after editing the sheet.
$ws.Range("H10").Select()
